# Refresh the crypto price/volume snapshot (Price = column D,
# Volume(1h) = column E) to match the latest scrape, row by row.
# Price cells are forced to Text format before the assignment so Excel
# doesn't silently coerce numeric-looking strings like "1.002" into
# floats/dates, then the cell style is reset back to "Normal" so no
# stray formatting is left behind (Volume cells already contain
# "%"/spaces so they stay text without any extra help).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.330.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4702"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2887"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.871.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.143"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6868"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "271.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.312.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007739"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.115.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.222"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.957"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09897"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.378"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.469"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04711"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7023"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01885"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.653"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.308"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4173"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.269"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.093"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "929.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05684"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
